$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.649.80"
$ws.Range("E2").Value = "  +3.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.699.56"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.40"
$ws.Range("E5").Value = "  +3.58%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4043"
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.543"
$ws.Range("E9").Value = "  +9.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.59"
$ws.Range("E10").Value = "  +11.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.000"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08830"
$ws.Range("E12").Value = "  +2.66%  "
$ws.Range("E13").Value = "  +8.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.47"
$ws.Range("E14").Value = "  +3.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001333"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.669"
$ws.Range("E16").Value = "  +6.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.701.73"
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "101.43"
$ws.Range("E18").Value = "  +1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07099"
$ws.Range("E19").Value = "  +4.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.83"
$ws.Range("E20").Value = "  +4.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.910"
$ws.Range("E21").Value = "  +4.29%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.18"
$ws.Range("E23").Value = "  +3.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.632.17"
$ws.Range("E24").Value = "  +3.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.172"
$ws.Range("E25").Value = "  +15.98%  "
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.48"
$ws.Range("E27").Value = "  +3.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.87"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.248"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.44"
$ws.Range("E30").Value = "  +3.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.622"
$ws.Range("E31").Value = "  +16.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.113"
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.883.81"
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.443"
$ws.Range("E34").Value = "  +13.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08596"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.51"
$ws.Range("E36").Value = "  +12.29%  "
$ws.Range("E37").Value = "  +5.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.945"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.82"
$ws.Range("E39").Value = "  +4.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02804"
$ws.Range("E40").Value = "  +11.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09122"
$ws.Range("E41").Value = "  +4.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7768"
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.465"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7296"
$ws.Range("E44").Value = "  +4.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.58"
$ws.Range("E45").Value = "  +4.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.518"
$ws.Range("E46").Value = "  +6.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.228"
$ws.Range("E47").Value = "  +4.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.386"
$ws.Range("E48").Value = "  +19.70%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.23"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08052"
$ws.Range("E51").Value = "  +4.28%  "
